$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.091.25'
$ws.Range("E2").Value = '  -4.85%  '
$ws.Range("D3").Value = '2.449.78'
$ws.Range("E3").Value = '  -3.96%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.88%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -3.41%  '
$ws.Range("D9").Value = '2.453.92'
$ws.Range("E9").Value = '  -4.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0989'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.27%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.90%  '
$ws.Range("E13").Value = '  -6.45%  '
$ws.Range("D14").Value = '2.886.84'
$ws.Range("E14").Value = '  -3.90%  '
$ws.Range("D15").Value = '57.983.18'
$ws.Range("E15").Value = '  -4.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000138'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.01%  '
$ws.Range("D18").Value = '2.455.12'
$ws.Range("E18").Value = '  -4.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.89%  '
$ws.Range("E20").Value = '  -5.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '317.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.60%  '
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.403'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.21%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.985'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("E28").Value = '  -7.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.94%  '
$ws.Range("D30").Value = '0.0₃0751'
$ws.Range("E30").Value = '  -7.43%  '
$ws.Range("E31").Value = '  -4.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("E34").Value = '  -13.43%  '
$ws.Range("E35").Value = '  -9.80%  '
$ws.Range("E36").Value = '  -4.31%  '
$ws.Range("E38").Value = '  -7.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = '  -6.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.777'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '270.26'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -12.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -13.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("E46").Value = '  -5.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0922'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.79%  '
$ws.Range("E49").Value = '  -5.16%  '
$ws.Range("E50").Value = '  -6.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.61%  '
